# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        13 = 2419
        16 = 10
        19 = 565
        24 = 1974
        25 = 4103
        30 = 2105
    }
    "全部类型" = @{
        13 = 2419
        17 = 10
        20 = 565
        25 = 1974
        26 = 4103
        31 = 2105
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
